# 0.4.0 bug fix, new enemy ashes-skull, eyeball of the elder done.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bug fix: SkeletonKnight's magical resistance (N17) was -0.25, now -0.2
$ws.Range("N17").Value = -0.2

# New enemy: insert a fresh row at row 18 (pushes everything below down by one)
$ws.Rows.Item(18).Insert() | Out-Null

$ws.Range("A18").Value = "AshesSkull"
$ws.Range("B18").Value = 20
$ws.Range("C18").Value = 1000
$ws.Range("D18").Value = 15
$ws.Range("E18").Value = 6
$ws.Range("F18").Value = 15
$ws.Range("I18").Value = 4
$ws.Range("J18").Value = 12
$ws.Range("K18").Value = "NORMAL"
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = 5
$ws.Range("N18").Value = 0.3
$ws.Range("O18").Value = "FIRE:0.3 SHADOW:0.2 HOLY:-0.5"
$ws.Range("P18").Value = "UNDEAD"

# Update the active selection to the new row
$ws.Range("N18").Select() | Out-Null
